$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 - this shifts existing rows 59:112 down to 60:113
# (matches the diff: dimension grows from A1:R112 to A1:R113, and every row
#  from the old 59 onward is now one row lower, with a brand-new record
#  appearing at row 59 and the old last record - row 112 - landing at row 113)
$ws.Rows(59).Insert()

# Populate the newly inserted row 59 with the new record's data.
$ws.Cells.Item(59, 1).Value = 4
$ws.Cells.Item(59, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value = "Los Lagos"
$ws.Cells.Item(59, 4).Value = 44634
$ws.Cells.Item(59, 5).Value = 10
$ws.Cells.Item(59, 6).Value = 100112022
$ws.Cells.Item(59, 7).Value = "Arveja Verde"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 20
$ws.Cells.Item(59, 11).Value = 30000
$ws.Cells.Item(59, 12).Value = 30000
$ws.Cells.Item(59, 13).Value = 30000
$ws.Cells.Item(59, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(59, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(59, 16).Value = 1200
$ws.Cells.Item(59, 17).Value = 25
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# Make sure the date cell keeps the workbook's date number format (style
# index 2 in the original file), same as every other "Fecha" cell in column D.
$ws.Cells.Item(59, 4).NumberFormat = $ws.Cells.Item(60, 4).NumberFormat
